# Apply the FlashScore weekly-games update for 2024-11-04:
#  - Refresh a batch of odds figures across rows 2-7
#  - Remove the Uruguay "Wanderers - Defensor Sp." fixture (row 8) entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ARGENTINA - TORNEO BETANO (Ind. Rivadavia - Rosario Central) ---
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 2.75
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.8
$ws.Range("P2").Value = 1.91
$ws.Range("Q2").Value = 3.6
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 2

# --- Row 3: BRAZIL - SERIE A BETANO (Corinthians - Palmeiras) ---
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 8
$ws.Range("AG3").Value = 7
$ws.Range("AU3").Value = 8.5
$ws.Range("AY3").Value = 26
$ws.Range("AZ3").Value = 51

# --- Row 4: BRAZIL - SERIE B (Ponte Preta - Paysandu PA) ---
$ws.Range("G4").Value = 2.3
$ws.Range("J4").Value = 3.1
$ws.Range("L4").Value = 4
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Z4").Value = 21
$ws.Range("AH4").Value = 15
$ws.Range("AO4").Value = 13
$ws.Range("AS4").Value = 251

# --- Row 5: BRAZIL - SERIE B (Ituano - CRB) ---
$ws.Range("G5").Value = 2.15
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7

# --- Row 6: COLOMBIA - PRIMERA A (Atl. Nacional - Santa Fe) ---
$ws.Range("G6").Value = 1.76
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63

# --- Row 7: PARAGUAY - PRIMERA DIVISION (Cerro Porteno - Tacuary) ---
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("Q7").Value = 1.9
$ws.Range("R7").Value = 1.9

# --- Remove row 8 (URUGUAY - PRIMERA DIVISION, Wanderers - Defensor Sp.) ---
$ws.Rows.Item(8).Delete()
